# Update "想去人数" (interested-count) figures on the 展览 and 全部类型
# sheets to match the freshly generated data snapshot.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value  = 4878
$ws1.Range("F8").Value  = 5
$ws1.Range("F9").Value  = 535
$ws1.Range("F10").Value = 489
$ws1.Range("F11").Value = 28
$ws1.Range("F13").Value = 1361
$ws1.Range("F14").Value = 3314
$ws1.Range("F15").Value = 394
$ws1.Range("F16").Value = 120
$ws1.Range("F17").Value = 105
$ws1.Range("F19").Value = 2522
$ws1.Range("F20").Value = 121
$ws1.Range("F23").Value = 173
$ws1.Range("F24").Value = 33
$ws1.Range("F27").Value = 253

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F6").Value  = 4878
$ws4.Range("F9").Value  = 5
$ws4.Range("F10").Value = 535
$ws4.Range("F11").Value = 489
$ws4.Range("F12").Value = 28
$ws4.Range("F14").Value = 1361
$ws4.Range("F15").Value = 3314
$ws4.Range("F16").Value = 394
$ws4.Range("F17").Value = 120
$ws4.Range("F18").Value = 105
$ws4.Range("F20").Value = 2522
$ws4.Range("F21").Value = 121
$ws4.Range("F24").Value = 173
$ws4.Range("F25").Value = 33
$ws4.Range("F28").Value = 253
